$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-12 down to 6-13
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with this week's data (latest record, pushed to the top of the data block)
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C5").Value = 'Ñuble'
$ws.Range("D5").Value = 44701
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112037
$ws.Range("G5").Value = 'Cebollín'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = 7250
$ws.Range("N5").Value = '$/paquete 36 unidades'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 201
$ws.Range("Q5").Value = 36
$ws.Range("R5").Value = 'Hortaliza'

# Ensure the date cell keeps the same date/time number format as the rest of column D
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
